$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C (Förändrad) for rows 2-28 from 45207 to 45208
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}

# 2) Update hyperlink formulas in row 2 (S2, T2, V2, W2, X2, Y2):
#    replace "Logging_SKELLEFTEA" with "Logging_2482" in the URL path
$cols = @("S", "T", "V", "W", "X", "Y")
foreach ($col in $cols) {
    $rangeCell = $ws.Range($col + "2")
    $formula = $rangeCell.Formula
    if ($formula) {
        $newFormula = $formula -replace "Logging_SKELLEFTEA", "Logging_2482"
        $rangeCell.Formula = $newFormula
    }
}
